$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update nombre_aides (col C) and montant_total (col D) for 2020-08-21 data refresh
$ws.Cells.Item(2, 3).Value = 38217
$ws.Cells.Item(2, 4).Value = 55264966
$ws.Cells.Item(3, 3).Value = 91916
$ws.Cells.Item(3, 4).Value = 134739319
$ws.Cells.Item(4, 3).Value = 31420
$ws.Cells.Item(4, 4).Value = 46533374
$ws.Cells.Item(5, 3).Value = 8790
$ws.Cells.Item(5, 4).Value = 13065563
$ws.Cells.Item(6, 3).Value = 2024
$ws.Cells.Item(6, 4).Value = 3007971
$ws.Cells.Item(7, 3).Value = 157
$ws.Cells.Item(7, 4).Value = 230593
$ws.Cells.Item(12, 3).Value = 41707
$ws.Cells.Item(12, 4).Value = 56584402
$ws.Cells.Item(13, 3).Value = 9781
$ws.Cells.Item(13, 4).Value = 14148708
$ws.Cells.Item(14, 3).Value = 26192
$ws.Cells.Item(14, 4).Value = 38407440
$ws.Cells.Item(15, 3).Value = 8374
$ws.Cells.Item(15, 4).Value = 12428177
$ws.Cells.Item(16, 3).Value = 2167
$ws.Cells.Item(16, 4).Value = 3222665
$ws.Cells.Item(20, 3).Value = 10316
$ws.Cells.Item(20, 4).Value = 13648192
$ws.Cells.Item(21, 3).Value = 13545
$ws.Cells.Item(21, 4).Value = 19550937
$ws.Cells.Item(22, 3).Value = 31900
$ws.Cells.Item(22, 4).Value = 46811435
$ws.Cells.Item(23, 3).Value = 10294
$ws.Cells.Item(23, 4).Value = 15302530
$ws.Cells.Item(24, 3).Value = 2662
$ws.Cells.Item(24, 4).Value = 3958182
$ws.Cells.Item(25, 3).Value = 512
$ws.Cells.Item(25, 4).Value = 762092
$ws.Cells.Item(27, 3).Value = 11779
$ws.Cells.Item(27, 4).Value = 15729043
$ws.Cells.Item(28, 3).Value = 7763
$ws.Cells.Item(28, 4).Value = 11237834
$ws.Cells.Item(29, 3).Value = 22732
$ws.Cells.Item(29, 4).Value = 33368834
$ws.Cells.Item(30, 3).Value = 7879
$ws.Cells.Item(30, 4).Value = 11722443
$ws.Cells.Item(31, 3).Value = 1991
$ws.Cells.Item(31, 4).Value = 2970999
$ws.Cells.Item(32, 3).Value = 369
$ws.Cells.Item(32, 4).Value = 550915
$ws.Cells.Item(34, 3).Value = 8385
$ws.Cells.Item(34, 4).Value = 11076411
$ws.Cells.Item(35, 3).Value = 3296
$ws.Cells.Item(35, 4).Value = 4759937
$ws.Cells.Item(36, 3).Value = 7933
$ws.Cells.Item(36, 4).Value = 11586137
$ws.Cells.Item(37, 3).Value = 3204
$ws.Cells.Item(37, 4).Value = 4748961
$ws.Cells.Item(38, 3).Value = 832
$ws.Cells.Item(38, 4).Value = 1239223
$ws.Cells.Item(39, 3).Value = 167
$ws.Cells.Item(39, 4).Value = 248186
$ws.Cells.Item(41, 3).Value = 2504
$ws.Cells.Item(41, 4).Value = 3386362
$ws.Cells.Item(42, 3).Value = 17459
$ws.Cells.Item(42, 4).Value = 25245824
$ws.Cells.Item(43, 3).Value = 51641
$ws.Cells.Item(43, 4).Value = 75700099
$ws.Cells.Item(44, 3).Value = 19136
$ws.Cells.Item(44, 4).Value = 28425008
$ws.Cells.Item(45, 3).Value = 5655
$ws.Cells.Item(45, 4).Value = 8418805
$ws.Cells.Item(46, 3).Value = 1222
$ws.Cells.Item(46, 4).Value = 1823545
$ws.Cells.Item(50, 3).Value = 16885
$ws.Cells.Item(50, 4).Value = 22457769
$ws.Cells.Item(51, 3).Value = 2083
$ws.Cells.Item(51, 4).Value = 3021046
$ws.Cells.Item(52, 3).Value = 7048
$ws.Cells.Item(52, 4).Value = 10359130
$ws.Cells.Item(53, 3).Value = 2388
$ws.Cells.Item(53, 4).Value = 3566464
$ws.Cells.Item(56, 3).Value = 20
$ws.Cells.Item(56, 4).Value = 30000
$ws.Cells.Item(57, 3).Value = 7152
$ws.Cells.Item(57, 4).Value = 9834800
$ws.Cells.Item(58, 3).Value = 1057
$ws.Cells.Item(58, 4).Value = 1695006
$ws.Cells.Item(59, 3).Value = 2643
$ws.Cells.Item(59, 4).Value = 4255513
$ws.Cells.Item(60, 3).Value = 1051
$ws.Cells.Item(60, 4).Value = 1704338
$ws.Cells.Item(61, 3).Value = 354
$ws.Cells.Item(61, 4).Value = 573883
$ws.Cells.Item(62, 3).Value = 117
$ws.Cells.Item(62, 4).Value = 194600
$ws.Cells.Item(63, 3).Value = 22
$ws.Cells.Item(63, 4).Value = 39000
$ws.Cells.Item(64, 3).Value = 1533
$ws.Cells.Item(64, 4).Value = 2298234
$ws.Cells.Item(65, 3).Value = 15574
$ws.Cells.Item(65, 4).Value = 22492082
$ws.Cells.Item(66, 3).Value = 45170
$ws.Cells.Item(66, 4).Value = 66095098
$ws.Cells.Item(67, 3).Value = 15830
$ws.Cells.Item(67, 4).Value = 23523052
$ws.Cells.Item(68, 3).Value = 4603
$ws.Cells.Item(68, 4).Value = 6855551
$ws.Cells.Item(69, 3).Value = 936
$ws.Cells.Item(69, 4).Value = 1392168
$ws.Cells.Item(73, 3).Value = 15226
$ws.Cells.Item(73, 4).Value = 20059511
$ws.Cells.Item(74, 3).Value = 52934
$ws.Cells.Item(74, 4).Value = 77036767
$ws.Cells.Item(75, 3).Value = 149222
$ws.Cells.Item(75, 4).Value = 219847817
$ws.Cells.Item(76, 3).Value = 64526
$ws.Cells.Item(76, 4).Value = 96152013
$ws.Cells.Item(77, 3).Value = 20652
$ws.Cells.Item(77, 4).Value = 30858322
$ws.Cells.Item(78, 3).Value = 4911
$ws.Cells.Item(78, 4).Value = 7335401
$ws.Cells.Item(85, 3).Value = 52234
$ws.Cells.Item(85, 4).Value = 71034478
$ws.Cells.Item(86, 3).Value = 4684
$ws.Cells.Item(86, 4).Value = 6788196
$ws.Cells.Item(87, 3).Value = 11720
$ws.Cells.Item(87, 4).Value = 17217020
$ws.Cells.Item(88, 3).Value = 3922
$ws.Cells.Item(88, 4).Value = 5845583
$ws.Cells.Item(89, 3).Value = 1357
$ws.Cells.Item(89, 4).Value = 2027989
$ws.Cells.Item(93, 3).Value = 5476
$ws.Cells.Item(93, 4).Value = 7362005
$ws.Cells.Item(94, 3).Value = 1619
$ws.Cells.Item(94, 4).Value = 2332033
$ws.Cells.Item(95, 3).Value = 5260
$ws.Cells.Item(95, 4).Value = 7748035
$ws.Cells.Item(96, 3).Value = 1959
$ws.Cells.Item(96, 4).Value = 2917476
$ws.Cells.Item(97, 3).Value = 700
$ws.Cells.Item(97, 4).Value = 1048960
$ws.Cells.Item(98, 3).Value = 189
$ws.Cells.Item(98, 4).Value = 282613
$ws.Cells.Item(101, 3).Value = 3617
$ws.Cells.Item(101, 4).Value = 4787798
$ws.Cells.Item(102, 3).Value = 670
$ws.Cells.Item(102, 4).Value = 1079209
$ws.Cells.Item(103, 3).Value = 401
$ws.Cells.Item(103, 4).Value = 656697
$ws.Cells.Item(104, 3).Value = 141
$ws.Cells.Item(104, 4).Value = 223160
$ws.Cells.Item(105, 3).Value = 50
$ws.Cells.Item(105, 4).Value = 81000
$ws.Cells.Item(107, 3).Value = 10921
$ws.Cells.Item(107, 4).Value = 15843648
$ws.Cells.Item(108, 3).Value = 29474
$ws.Cells.Item(108, 4).Value = 43299276
$ws.Cells.Item(109, 3).Value = 9864
$ws.Cells.Item(109, 4).Value = 14668389
$ws.Cells.Item(110, 3).Value = 2716
$ws.Cells.Item(110, 4).Value = 4050207
$ws.Cells.Item(114, 3).Value = 9877
$ws.Cells.Item(114, 4).Value = 13045828
$ws.Cells.Item(115, 3).Value = 30833
$ws.Cells.Item(115, 4).Value = 44461536
$ws.Cells.Item(116, 3).Value = 66761
$ws.Cells.Item(116, 4).Value = 97698181
$ws.Cells.Item(117, 3).Value = 21542
$ws.Cells.Item(117, 4).Value = 32015213
$ws.Cells.Item(118, 3).Value = 6111
$ws.Cells.Item(118, 4).Value = 9104521
$ws.Cells.Item(119, 3).Value = 1140
$ws.Cells.Item(119, 4).Value = 1703771
$ws.Cells.Item(124, 3).Value = 26095
$ws.Cells.Item(124, 4).Value = 34841535
$ws.Cells.Item(125, 3).Value = 36475
$ws.Cells.Item(125, 4).Value = 52638126
$ws.Cells.Item(126, 3).Value = 77553
$ws.Cells.Item(126, 4).Value = 113401904
$ws.Cells.Item(127, 3).Value = 24042
$ws.Cells.Item(127, 4).Value = 35683287
$ws.Cells.Item(128, 3).Value = 6441
$ws.Cells.Item(128, 4).Value = 9572238
$ws.Cells.Item(129, 3).Value = 1255
$ws.Cells.Item(129, 4).Value = 1865811
$ws.Cells.Item(133, 3).Value = 32083
$ws.Cells.Item(133, 4).Value = 42598137
$ws.Cells.Item(134, 3).Value = 13425
$ws.Cells.Item(134, 4).Value = 19436539
$ws.Cells.Item(135, 3).Value = 32621
$ws.Cells.Item(135, 4).Value = 47909991
$ws.Cells.Item(136, 3).Value = 11556
$ws.Cells.Item(136, 4).Value = 17169433
$ws.Cells.Item(137, 3).Value = 2985
$ws.Cells.Item(137, 4).Value = 4448741
$ws.Cells.Item(138, 3).Value = 506
$ws.Cells.Item(138, 4).Value = 752990
$ws.Cells.Item(141, 3).Value = 10907
$ws.Cells.Item(141, 4).Value = 14542812
$ws.Cells.Item(142, 3).Value = 35571
$ws.Cells.Item(142, 4).Value = 51376519
$ws.Cells.Item(143, 3).Value = 82195
$ws.Cells.Item(143, 4).Value = 120422577
$ws.Cells.Item(144, 3).Value = 24576
$ws.Cells.Item(144, 4).Value = 36512527
$ws.Cells.Item(145, 3).Value = 6454
$ws.Cells.Item(145, 4).Value = 9630567
$ws.Cells.Item(146, 3).Value = 1459
$ws.Cells.Item(146, 4).Value = 2170730
$ws.Cells.Item(149, 3).Value = 29494
$ws.Cells.Item(149, 4).Value = 39777057
